$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 1: C1 16 -> 30, D1 FALSE -> TRUE
$ws.Range("C1").Value = 30
$ws.Range("D1").Value = $true

# Row 2: C2 1 -> 5, D2 0 -> 50, E2 FALSE -> TRUE
$ws.Range("C2").Value = 5
$ws.Range("D2").Value = 50
$ws.Range("E2").Value = $true

# Row 3: C3 100 -> 52, D3 stays TRUE (unchanged)
$ws.Range("C3").Value = 52
